# Apply updated dSF (column F) values after re-pulling data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -6
$ws.Range("F11").Value = -2
$ws.Range("F14").Value = 0
